$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Update period dates (B) and (C) for rows 8-10 (date serials, no time component)
$ws.Range("B8").Value2 = 44835
$ws.Range("C8").Value2 = 44926
$ws.Range("B9").Value2 = 44835
$ws.Range("C9").Value2 = 44926
$ws.Range("B10").Value2 = 44835
$ws.Range("C10").Value2 = 44926

# Update validation (AA) and update (AB) dates for rows 8-10
$ws.Range("AA8").Value2 = 44936
$ws.Range("AB8").Value2 = 44936
$ws.Range("AA9").Value2 = 44936
$ws.Range("AB9").Value2 = 44936
$ws.Range("AA10").Value2 = 44936
$ws.Range("AB10").Value2 = 44936

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:C2").Select()
